$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.495.35'
$ws.Range('E2').Value = '  +2.02%  '
$ws.Range('D3').Value = '1.856.16'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6962'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.93%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07692'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3072'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.60'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07783'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = '1.867.30'
$ws.Range('E12').Value = '  +2.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.161'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.61%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '91.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6934'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.294'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.29%  '
$ws.Range('D17').Value = '29.478.12'
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008339'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '2.098.55'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '238.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.74'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9992'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.611'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9995'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.67%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.887'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.531'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.246'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.150'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.199'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.30%  '
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7745'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.885'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.150'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.688'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').Value = '1.315.09'
$ws.Range('E38').Value = '  +7.63%  '
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.722'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9516'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.14'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.781'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9995'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.823'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.68%  '
$ws.Range('E46').Value = '  +1.89%  '
$ws.Range('D47').Value = '2.001.69'
$ws.Range('E47').Value = '  +1.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5217'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.788'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '63.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.958'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.89%  '
